# "contingencies with rene fine"
#
# Two new line rows (line7, line8) are inserted into the lines table right
# after line6 (row 7), pushing the existing extr1..extr8 rows (previously
# rows 8-15) down by two rows (to rows 10-17). Numeric values across the
# whole table (C/D/E columns) are also refreshed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: make room for two new rows right after "line6" (row 7), before
# the "extr1" row (row 8).
$ws.Rows.Item(8).Resize(2).Insert(-4121) | Out-Null   # xlShiftDown

# Column A in this table uses the bold/bordered/centered "name" style (same
# style as the other row-label cells); reapply it to the two new A cells so
# they match their neighbours instead of the blank default style the insert
# leaves behind.
$ws.Cells.Item(7, 1).Copy() | Out-Null
$ws.Range($ws.Cells.Item(8, 1), $ws.Cells.Item(9, 1)).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Step 2: populate the two new rows with the line7 / line8 data.
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# --- Step 3: refresh the rest of the table (formerly extr1..extr8, now
# living two rows further down at rows 10-17). The name text itself
# (column B) already moved down correctly with the row insert. Column A is
# a plain running index (0-based), so it needs to be bumped by 2 to stay
# sequential now that two rows were inserted above it; columns C/D/E get
# their refreshed numeric/boolean values.
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = $false

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $true

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $false

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $true

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $false

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $false
